$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.828516
$ws.Range("H2").Value = 2.485548
$ws.Range("I2").Value = 0.4625620436231038
$ws.Range("J2").Value = 0.4821955800271095
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 16.60495066666667
$ws.Range("N2").Value = 49.814852
$ws.Range("O2").Value = 0.4330603147186406
$ws.Range("P2").Value = 0.5197056776409935
$ws.Range("Q2").Value = 13.757467306544
$ws.Range("R2").Value = 123.817205758896
$ws.Range("S2").Value = 0.2003172641883189
$ws.Range("T2").Value = 0.2505997806734809

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.828516
$ws.Range("H3").Value = 2.485548
$ws.Range("I3").Value = 0.4625620436231038
$ws.Range("J3").Value = 0.4821955800271095
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.7472513333333333
$ws.Range("N3").Value = 2.241754
$ws.Range("O3").Value = 0.01948845883877707
$ws.Range("P3").Value = 0.02338764916283215
$ws.Range("Q3").Value = 0.619109685688
$ws.Range("R3").Value = 5.571987171192
$ws.Range("S3").Value = 0.009014621347529461
$ws.Range("T3").Value = 0.01127742105354239

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.828516
$ws.Range("H4").Value = 2.485548
$ws.Range("I4").Value = 0.4625620436231038
$ws.Range("J4").Value = 0.4821955800271095
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.5008306666666666
$ws.Range("N4").Value = 1.502492
$ws.Range("O4").Value = 0.01306176034372721
$ws.Range("P4").Value = 0.01567511679067463
$ws.Range("Q4").Value = 0.414946220624
$ws.Range("R4").Value = 3.734515985616
$ws.Range("S4").Value = 0.00604187455790967
$ws.Range("T4").Value = 0.007558472032872035

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.828516
$ws.Range("H5").Value = 2.485548
$ws.Range("I5").Value = 0.4625620436231038
$ws.Range("J5").Value = 0.4821955800271095
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.312462666666667
$ws.Range("N5").Value = 3.937388
$ws.Range("O5").Value = 0.034229279381366
$ws.Range("P5").Value = 0.04107776730272161
$ws.Range("Q5").Value = 1.087396318736
$ws.Range("R5").Value = 9.786566868624
$ws.Range("S5").Value = 0.01583316542239083
$ws.Range("T5").Value = 0.01980751783075448

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.828516
$ws.Range("H6").Value = 2.485548
$ws.Range("I6").Value = 0.4625620436231038
$ws.Range("J6").Value = 0.4821955800271095
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 19.1777795
$ws.Range("N6").Value = 38.355559
$ws.Range("O6").Value = 0.5001601867174891
$ws.Range("P6").Value = 0.4001537891027781
$ws.Range("Q6").Value = 15.889097160222
$ws.Range("R6").Value = 95.334582961332
$ws.Range("S6").Value = 0.2313551181069549
$ws.Range("T6").Value = 0.1929523884364598

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.7438396666666667
$ws.Range("H7").Value = 2.231519
$ws.Range("I7").Value = 0.4152870872032183
$ws.Range("J7").Value = 0.4329140288365043
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 16.60495066666667
$ws.Range("N7").Value = 49.814852
$ws.Range("O7").Value = 0.4330603147186406
$ws.Range("P7").Value = 0.5197056776409935
$ws.Range("Q7").Value = 12.35142096890978
$ws.Range("R7").Value = 111.162788720188
$ws.Range("S7").Value = 0.1798443566828133
$ws.Range("T7").Value = 0.2249878787167681

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.7438396666666667
$ws.Range("H8").Value = 2.231519
$ws.Range("I8").Value = 0.4152870872032183
$ws.Range("J8").Value = 0.4329140288365043
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.7472513333333333
$ws.Range("N8").Value = 2.241754
$ws.Range("O8").Value = 0.01948845883877707
$ws.Range("P8").Value = 0.02338764916283215
$ws.Range("Q8").Value = 0.5558351827028889
$ws.Range("R8").Value = 5.002516644326
$ws.Range("S8").Value = 0.008093305305235544
$ws.Range("T8").Value = 0.01012484142409636

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.7438396666666667
$ws.Range("H9").Value = 2.231519
$ws.Range("I9").Value = 0.4152870872032183
$ws.Range("J9").Value = 0.4329140288365043
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.5008306666666666
$ws.Range("N9").Value = 1.502492
$ws.Range("O9").Value = 0.01306176034372721
$ws.Range("P9").Value = 0.01567511679067463
$ws.Range("Q9").Value = 0.3725377161497778
$ws.Range("R9").Value = 3.352839445348
$ws.Range("S9").Value = 0.005424380406892979
$ws.Range("T9").Value = 0.006785977962333687

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.7438396666666667
$ws.Range("H10").Value = 2.231519
$ws.Range("I10").Value = 0.4152870872032183
$ws.Range("J10").Value = 0.4329140288365043
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.312462666666667
$ws.Range("N10").Value = 3.937388
$ws.Range("O10").Value = 0.034229279381366
$ws.Range("P10").Value = 0.04107776730272161
$ws.Range("Q10").Value = 0.9762617924857777
$ws.Range("R10").Value = 8.786356132371999
$ws.Range("S10").Value = 0.01421497773135267
$ws.Range("T10").Value = 0.01778314173862963

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.7438396666666667
$ws.Range("H11").Value = 2.231519
$ws.Range("I11").Value = 0.4152870872032183
$ws.Range("J11").Value = 0.4329140288365043
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = 19.1777795
$ws.Range("N11").Value = 38.355559
$ws.Range("O11").Value = 0.5001601867174891
$ws.Range("P11").Value = 0.4001537891027781
$ws.Range("Q11").Value = 14.26519311068683
$ws.Range("R11").Value = 85.591158664121
$ws.Range("S11").Value = 0.2077100670769239
$ws.Range("T11").Value = 0.1732321889946765

$ws.Range("E12").Value = 2
$ws.Range("G12").Value = 0.21879
$ws.Range("H12").Value = 0.43758
$ws.Range("I12").Value = 0.1221508691736778
$ws.Range("J12").Value = 0.08489039113638626
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 16.60495066666667
$ws.Range("N12").Value = 49.814852
$ws.Range("O12").Value = 0.4330603147186406
$ws.Range("P12").Value = 0.5197056776409935
$ws.Range("Q12").Value = 3.63299715636
$ws.Range("R12").Value = 21.79798293816
$ws.Range("S12").Value = 0.05289869384750843
$ws.Range("T12").Value = 0.04411801825074461

$ws.Range("E13").Value = 2
$ws.Range("G13").Value = 0.21879
$ws.Range("H13").Value = 0.43758
$ws.Range("I13").Value = 0.1221508691736778
$ws.Range("J13").Value = 0.08489039113638626
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.7472513333333333
$ws.Range("N13").Value = 2.241754
$ws.Range("O13").Value = 0.01948845883877707
$ws.Range("P13").Value = 0.02338764916283215
$ws.Range("Q13").Value = 0.16349111922
$ws.Range("R13").Value = 0.9809467153199999
$ws.Range("S13").Value = 0.002380532186012063
$ws.Range("T13").Value = 0.001985386685193398

$ws.Range("E14").Value = 2
$ws.Range("G14").Value = 0.21879
$ws.Range("H14").Value = 0.43758
$ws.Range("I14").Value = 0.1221508691736778
$ws.Range("J14").Value = 0.08489039113638626
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 0.5008306666666666
$ws.Range("N14").Value = 1.502492
$ws.Range("O14").Value = 0.01306176034372721
$ws.Range("P14").Value = 0.01567511679067463
$ws.Range("Q14").Value = 0.10957674156
$ws.Range("R14").Value = 0.6574604493599999
$ws.Range("S14").Value = 0.001595505378924555
$ws.Range("T14").Value = 0.001330666795468905

$ws.Range("E15").Value = 2
$ws.Range("G15").Value = 0.21879
$ws.Range("H15").Value = 0.43758
$ws.Range("I15").Value = 0.1221508691736778
$ws.Range("J15").Value = 0.08489039113638626
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 1.312462666666667
$ws.Range("N15").Value = 3.937388
$ws.Range("O15").Value = 0.034229279381366
$ws.Range("P15").Value = 0.04107776730272161
$ws.Range("Q15").Value = 0.28715370684
$ws.Range("R15").Value = 1.72292224104
$ws.Range("S15").Value = 0.004181136227622507
$ws.Range("T15").Value = 0.003487107733337495

$ws.Range("E16").Value = 2
$ws.Range("G16").Value = 0.21879
$ws.Range("H16").Value = 0.43758
$ws.Range("I16").Value = 0.1221508691736778
$ws.Range("J16").Value = 0.08489039113638626
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 19.1777795
$ws.Range("N16").Value = 38.355559
$ws.Range("O16").Value = 0.5001601867174891
$ws.Range("P16").Value = 0.4001537891027781
$ws.Range("Q16").Value = 4.195906376805
$ws.Range("R16").Value = 16.78362550722
$ws.Range("S16").Value = 0.06109500153361029
$ws.Range("T16").Value = 0.03396921167164185
